# Update cryptocurrency price/volume data as scraped on the latest run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "64.415.20"
Set-TextValue $ws.Range("E2") "  +0.17%  "
Set-TextValue $ws.Range("D3") "3.506.67"
Set-TextValue $ws.Range("E3") "  +0.57%  "
Set-TextValue $ws.Range("E4") "  +0.04%  "
Set-TextValue $ws.Range("D5") "591.28"
Set-TextValue $ws.Range("E5") "  +0.87%  "
Set-TextValue $ws.Range("D6") "134.62"
Set-TextValue $ws.Range("E6") "  +0.28%  "
Set-TextValue $ws.Range("E7") "  +0.01%  "
Set-TextValue $ws.Range("E8") "  +0.70%  "
Set-TextValue $ws.Range("D9") "7.64"
Set-TextValue $ws.Range("E9") "  +6.06%  "
Set-TextValue $ws.Range("E10") "  +1.49%  "
Set-TextValue $ws.Range("D11") "0.391"
Set-TextValue $ws.Range("E11") "  +4.18%  "
Set-TextValue $ws.Range("E13") "  +0.68%  "
Set-TextValue $ws.Range("E14") "  +1.20%  "
Set-TextValue $ws.Range("D15") "3.506.46"
Set-TextValue $ws.Range("E15") "  +0.61%  "
Set-TextValue $ws.Range("D16") "25.83"
Set-TextValue $ws.Range("E16") "  +2.92%  "
Set-TextValue $ws.Range("D17") "64.409.96"
Set-TextValue $ws.Range("E17") "  +0.14%  "
Set-TextValue $ws.Range("D18") "10.06"
Set-TextValue $ws.Range("E18") "  +0.67%  "
Set-TextValue $ws.Range("D19") "5.80"
Set-TextValue $ws.Range("E19") "  +1.95%  "
Set-TextValue $ws.Range("E20") "  -0.29%  "
Set-TextValue $ws.Range("D21") "391.85"
Set-TextValue $ws.Range("E21") "  +2.04%  "
Set-TextValue $ws.Range("D22") "0.584"
Set-TextValue $ws.Range("E22") "  +3.38%  "
Set-TextValue $ws.Range("D23") "3.646.45"
Set-TextValue $ws.Range("E23") "  +0.55%  "
Set-TextValue $ws.Range("D24") "74.49"
Set-TextValue $ws.Range("E24") "  +0.40%  "
Set-TextValue $ws.Range("E25") "  -0.07%  "
Set-TextValue $ws.Range("D26") "5.66"
Set-TextValue $ws.Range("E26") "  -0.49%  "
Set-TextValue $ws.Range("E27") "  +4.51%  "
Set-TextValue $ws.Range("E28") "  +0.85%  "
Set-TextValue $ws.Range("D29") "0.999"
Set-TextValue $ws.Range("E29") "  -0.03%  "
Set-TextValue $ws.Range("D30") "2.27"
Set-TextValue $ws.Range("E30") "  +2.01%  "
Set-TextValue $ws.Range("D31") "8.22"
Set-TextValue $ws.Range("E31") "  +0.05%  "
Set-TextValue $ws.Range("E32") "  -4.36%  "
Set-TextValue $ws.Range("D33") "0.157"
Set-TextValue $ws.Range("E33") "  +7.43%  "
Set-TextValue $ws.Range("D34") "3.534.98"
Set-TextValue $ws.Range("E34") "  +0.70%  "
Set-TextValue $ws.Range("E35") "  +0.01%  "
Set-TextValue $ws.Range("E36") "  +0.38%  "
Set-TextValue $ws.Range("D37") "5.36"
Set-TextValue $ws.Range("E37") "  +2.12%  "
Set-TextValue $ws.Range("E38") "  +2.08%  "
Set-TextValue $ws.Range("D39") "1.57"
Set-TextValue $ws.Range("E39") "  +2.71%  "
Set-TextValue $ws.Range("D40") "165.43"
Set-TextValue $ws.Range("E40") "  +2.08%  "
Set-TextValue $ws.Range("D41") "0.0793"
Set-TextValue $ws.Range("E41") "  +2.17%  "
Set-TextValue $ws.Range("E42") "  +0.84%  "
Set-TextValue $ws.Range("E43") "  +0.06%  "
Set-TextValue $ws.Range("D44") "4.46"
Set-TextValue $ws.Range("E44") "  +1.72%  "
Set-TextValue $ws.Range("D45") "24.95"
Set-TextValue $ws.Range("E45") "  -1.95%  "
Set-TextValue $ws.Range("E46") "  -0.84%  "
Set-TextValue $ws.Range("D47") "1.66"
Set-TextValue $ws.Range("E47") "  +1.48%  "
Set-TextValue $ws.Range("D48") "0.930"
Set-TextValue $ws.Range("E48") "  +4.08%  "
Set-TextValue $ws.Range("D49") "6.83"
Set-TextValue $ws.Range("E49") "  +1.62%  "
Set-TextValue $ws.Range("D50") "2.419.82"
Set-TextValue $ws.Range("E50") "  -1.83%  "
Set-TextValue $ws.Range("E51") "  +0.69%  "
